$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark rows 25, 26, 27, 30, 31 as achievements with an "x" flag in column A,
# matching the other already-implemented achievement rows.
$ws.Range("A25").Value = "x"
$ws.Range("A26").Value = "x"
$ws.Range("A27").Value = "x"
$ws.Range("A30").Value = "x"
$ws.Range("A31").Value = "x"

# Move the active selection to B22 (was B27).
$ws.Range("B22").Select()
